$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.508.85"
$ws.Range("E2").Value = "  +5.16%  "
$ws.Range("D3").Value = "2.231.30"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'228.84"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("D7").Value = "'61.77"
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("D9").Value = "'0.403"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").Value = "'58.41"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "'0.0877"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "2.564.78"
$ws.Range("E13").Value = "  +3.55%  "
$ws.Range("D14").Value = "'15.57"
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "'21.84"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "'5.58"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").Value = "2.238.71"
$ws.Range("E18").Value = "  +3.45%  "
$ws.Range("D19").Value = "41.457.14"
$ws.Range("E19").Value = "  +4.58%  "
$ws.Range("D20").Value = "'73.32"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  +6.13%  "
$ws.Range("D22").Value = "'6.11"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "'247.00"
$ws.Range("E23").Value = "  +7.63%  "
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").Value = "'9.53"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'168.91"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.142"
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("D30").Value = "'20.02"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "'1.44"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").Value = "'2.79"
$ws.Range("E32").Value = "  +4.58%  "
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "'4.97"
$ws.Range("E34").Value = "  +5.79%  "
$ws.Range("D35").Value = "'4.61"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'0.0623"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").Value = "'3.79"
$ws.Range("E37").Value = "  +4.76%  "
$ws.Range("D38").Value = "'6.66"
$ws.Range("E38").Value = "  -4.62%  "
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "'0.000233"
$ws.Range("E41").Value = "  +20.91%  "
$ws.Range("D42").Value = "'4.80"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("D43").Value = "'8.81"
$ws.Range("E43").Value = "  +14.20%  "
$ws.Range("E44").Value = "  +3.88%  "
$ws.Range("D45").Value = "'99.64"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").Value = "'0.0962"
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("D47").Value = "1.485.19"
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("D48").Value = "'1.19"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Value = "'16.37"
$ws.Range("E49").Value = "  -8.11%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  -1.71%  "
